$d = $word.ActiveDocument

$find = "быстро , гордитесь"
$replace = "быстро ,покажите анимации которые я добавил пожалуйста, пусть посмотрят на то как поднимается опускается блок со входом, и на остальные анимации, сконцентрируйте их внимание на этом, это очень важно, это даст нам больше шансов на победу гордитесь"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
